$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "24.667.26"
$ws.Cells.Item(2, 5).Value = "  +0.13%  "

# Row 3
Set-TextValue 3 4 "1.702.80"
$ws.Cells.Item(3, 5).Value = "  +0.86%  "

# Row 4
Set-TextValue 4 4 "0.9981"
$ws.Cells.Item(4, 5).Value = "  -0.76%  "

# Row 5
Set-TextValue 5 4 "314.02"
$ws.Cells.Item(5, 5).Value = "  -0.45%  "

# Row 6
Set-TextValue 6 4 "0.9966"
$ws.Cells.Item(6, 5).Value = "  -0.83%  "

# Row 7
Set-TextValue 7 4 "0.3976"
$ws.Cells.Item(7, 5).Value = "  +0.42%  "

# Row 8
Set-TextValue 8 4 "0.4069"
$ws.Cells.Item(8, 5).Value = "  +1.60%  "

# Row 9
$ws.Cells.Item(9, 2).Value = "BinanceUSD"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue 9 4 "0.9954"
$ws.Cells.Item(9, 5).Value = "  -1.03%  "

# Row 10
$ws.Cells.Item(10, 2).Value = "Polygon"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue 10 4 "1.512"
$ws.Cells.Item(10, 5).Value = "  +6.39%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +9.18%  "

# Row 12
Set-TextValue 12 4 "0.08792"
$ws.Cells.Item(12, 5).Value = "  -0.08%  "

# Row 13
Set-TextValue 13 4 "7.349"
$ws.Cells.Item(13, 5).Value = "  +11.04%  "

# Row 14
Set-TextValue 14 4 "23.35"
$ws.Cells.Item(14, 5).Value = "  +0.30%  "

# Row 15
Set-TextValue 15 4 "0.00001322"
$ws.Cells.Item(15, 5).Value = "  -0.25%  "

# Row 16
Set-TextValue 16 4 "7.530"
$ws.Cells.Item(16, 5).Value = "  +3.86%  "

# Row 17
Set-TextValue 17 4 "1.703.43"
$ws.Cells.Item(17, 5).Value = "  +0.41%  "

# Row 18
Set-TextValue 18 4 "100.97"
$ws.Cells.Item(18, 5).Value = "  -1.42%  "

# Row 19
Set-TextValue 19 4 "0.07101"
$ws.Cells.Item(19, 5).Value = "  +3.66%  "

# Row 20
Set-TextValue 20 4 "19.53"
$ws.Cells.Item(20, 5).Value = "  -0.70%  "

# Row 21
Set-TextValue 21 4 "6.765"
$ws.Cells.Item(21, 5).Value = "  -0.81%  "

# Row 22
Set-TextValue 22 4 "0.9963"
$ws.Cells.Item(22, 5).Value = "  -0.76%  "

# Row 23
Set-TextValue 23 4 "14.22"
$ws.Cells.Item(23, 5).Value = "  +1.54%  "

# Row 24
Set-TextValue 24 4 "24.717.06"
$ws.Cells.Item(24, 5).Value = "  +0.38%  "

# Row 25
Set-TextValue 25 4 "3.007"
$ws.Cells.Item(25, 5).Value = "  +5.89%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.22%  "

# Row 27
Set-TextValue 27 4 "22.44"
$ws.Cells.Item(27, 5).Value = "  +0.86%  "

# Row 28
Set-TextValue 28 4 "159.23"
$ws.Cells.Item(28, 5).Value = "  -0.45%  "

# Row 29
Set-TextValue 29 4 "5.117"
$ws.Cells.Item(29, 5).Value = "  -3.37%  "

# Row 30
Set-TextValue 30 4 "133.50"
$ws.Cells.Item(30, 5).Value = "  -0.13%  "

# Row 31
Set-TextValue 31 4 "7.447"
$ws.Cells.Item(31, 5).Value = "  +25.63%  "

# Row 32
Set-TextValue 32 4 "1.887.81"
$ws.Cells.Item(32, 5).Value = "  +0.10%  "

# Row 33
Set-TextValue 33 4 "1.087"
$ws.Cells.Item(33, 5).Value = "  -8.47%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 34 4 "7.433"
$ws.Cells.Item(34, 5).Value = "  +20.02%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Hedera"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 35 4 "0.08680"
$ws.Cells.Item(35, 5).Value = "  -2.39%  "

# Row 36
Set-TextValue 36 4 "11.10"
$ws.Cells.Item(36, 5).Value = "  +0.83%  "

# Row 37
Set-TextValue 37 4 "1.954"
$ws.Cells.Item(37, 5).Value = "  +3.87%  "

# Row 38
Set-TextValue 38 4 "0.2739"
$ws.Cells.Item(38, 5).Value = "  +0.94%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -3.79%  "

# Row 40
Set-TextValue 40 4 "0.02798"
$ws.Cells.Item(40, 5).Value = "  +10.15%  "

# Row 41
Set-TextValue 41 4 "0.09006"
$ws.Cells.Item(41, 5).Value = "  +0.96%  "

# Row 42
Set-TextValue 42 4 "1.484"
$ws.Cells.Item(42, 5).Value = "  +1.38%  "

# Row 43
Set-TextValue 43 4 "0.7674"
$ws.Cells.Item(43, 5).Value = "  +0.53%  "

# Row 44
Set-TextValue 44 4 "0.7223"
$ws.Cells.Item(44, 5).Value = "  +0.67%  "

# Row 45
Set-TextValue 45 4 "15.55"
$ws.Cells.Item(45, 5).Value = "  +1.06%  "

# Row 46
Set-TextValue 46 4 "2.464"
$ws.Cells.Item(46, 5).Value = "  -0.22%  "

# Row 47
Set-TextValue 47 4 "4.165"
$ws.Cells.Item(47, 5).Value = "  +0.70%  "

# Row 48
Set-TextValue 48 4 "0.9962"
$ws.Cells.Item(48, 5).Value = "  -0.78%  "

# Row 49
Set-TextValue 49 4 "141.59"
$ws.Cells.Item(49, 5).Value = "  -0.50%  "

# Row 50
Set-TextValue 50 4 "1.320"
$ws.Cells.Item(50, 5).Value = "  +13.46%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue 51 4 "0.00000000375"
$ws.Cells.Item(51, 5).Value = "  -3.52%  "
